$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Add new "type" column (D) with header and sample value
$ws.Range("D1").Value = "type"
$ws.Range("D2").Value = "1(deposit), 2(withdraw)"

# Move the active selection to match the committed workbook state
$ws.Range("F8").Select()
